$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8.694299999999998
$ws.Range("A9").Value = -22.2031
$ws.Range("B9").Value = 6.436599999999999
$ws.Range("B11").Value = 4.877600000000005
$ws.Range("A13").Value = -22.28530000000002
$ws.Range("A16").Value = -21.6193
$ws.Range("B16").Value = 4.782699999999998
$ws.Range("A18").Value = -22.13800000000002
$ws.Range("A20").Value = -21.47419999999998
$ws.Range("B23").Value = 8.928599999999996
$ws.Range("B24").Value = 4.883999999999999
$ws.Range("A26").Value = -21.09339999999996
$ws.Range("B26").Value = 5.163400000000003
$ws.Range("A27").Value = -21.77929999999999
$ws.Range("A29").Value = -21.70099999999997
$ws.Range("B34").Value = 9.462700000000007
$ws.Range("A35").Value = -21.91389999999999
$ws.Range("B35").Value = 5.481900000000003
$ws.Range("A36").Value = -20.80439999999998
$ws.Range("B44").Value = 4.431200000000005
$ws.Range("A45").Value = -21.42199999999998
$ws.Range("B48").Value = 5.780600000000001
$ws.Range("B49").Value = 5.7537
$ws.Range("B52").Value = 5.846399999999997
$ws.Range("A55").Value = -22.196
$ws.Range("A57").Value = -22.31980000000001
$ws.Range("B66").Value = 4.954399999999996
$ws.Range("B67").Value = 5.247399999999999
$ws.Range("A69").Value = -21.65779999999997
$ws.Range("B73").Value = 9.236399999999994
$ws.Range("A76").Value = -19.41279999999998
$ws.Range("A78").Value = -21.8184
$ws.Range("B78").Value = 5.882999999999999
$ws.Range("B80").Value = 9.527799999999996
$ws.Range("A82").Value = -21.68239999999999
$ws.Range("A83").Value = -21.67669999999999
$ws.Range("B91").Value = 5.084799999999998
$ws.Range("A93").Value = -21.48120000000002
$ws.Range("A97").Value = -21.62599999999998
$ws.Range("B97").Value = 5.022599999999995
$ws.Range("B99").Value = 5.357500000000001
$ws.Range("B104").Value = 9.940500000000004
